$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'erben',
    'reizen',
    'wenden',
    'lehnen',
    'fallen',
    'ruhen',
    'wüten',
    'schenken',
    'knurren',
    'knien',
    'gründen',
    'kichern',
    'klettern',
    'suchen',
    'parken',
    'feiern',
    'achten',
    'räumen',
    'gnaden',
    'wundern',
    'treiben',
    'gelten',
    'bluten',
    'kümmern',
    'bauen',
    'schicken',
    'doppeln',
    'hauen',
    'wiegen',
    'schrecken',
    'sichern',
    'wärmen',
    'starren',
    'irren',
    'werfen',
    'seufzen',
    'warten',
    'pflegen',
    'orten',
    'brauchen',
    'runden',
    'drohen',
    'malen',
    'filmen',
    'formen',
    'lohnen',
    'sprengen',
    'schneiden',
    'schwingen',
    'wählen',
    'saufen',
    'zögern',
    'buchen',
    'nennen',
    'schwören',
    'sperren',
    'schwächen',
    'kaufen',
    'platzen',
    'mögen',
    'schwimmen',
    'münzen',
    'wahren',
    'füttern',
    'sterben',
    'klingen',
    'tropfen',
    'lesen',
    'wollen',
    'schlucken',
    'wirken',
    'folgen',
    'spinnen',
    'bergen',
    'gleichen',
    'machen',
    'beißen',
    'wachsen',
    'stellen',
    'schmecken',
    'stecken',
    'trauen',
    'kratzen',
    'wehtun',
    'rühren',
    'sorgen',
    'baden',
    'zählen',
    'spielen',
    'planen',
    'quälen',
    'flehen',
    'führen',
    'klappen',
    'pflanzen',
    'kosten',
    'altern',
    'heilen',
    'arten',
    'plaudern',
    'boxen',
    'lügen',
    'jubeln',
    'betteln',
    'hören',
    'stammen',
    'lenken',
    'lösen',
    'treffen',
    'schreiten',
    'läuten',
    'duschen',
    'knarren',
    'backen',
    'loben',
    'schämen',
    'fließen',
    'fahren',
    'dringen',
    'schütteln',
    'nutzen',
    'liefern',
    'fischen',
    'grüßen',
    'atmen',
    'helfen',
    'liegen',
    'flüchten',
    'tollen',
    'stoßen',
    'warnen',
    'deuten',
    'geben',
    'lockern',
    'ärgern',
    'rufen',
    'fällen',
    'weinen',
    'heulen',
    'dienen',
    'bitten',
    'graben',
    'reiten',
    'stehlen',
    'töten',
    'hacken',
    'schulden',
    'tauschen',
    'äußern',
    'sinken',
    'scheitern',
    'weichen',
    'kehren',
    'sagen',
    'pfeifen',
    'zünden',
    'biegen',
    'fangen',
    'enden',
    'bremsen',
    'herrschen',
    'wagen',
    'ehren',
    'streichen',
    'spüren',
    'zeigen',
    'scheinen',
    'nähen',
    'stören',
    'mauern',
    'trennen',
    'rasen',
    'ändern',
    'werden',
    'bellen',
    'klingeln',
    'jagen',
    'zeichnen',
    'eignen',
    'decken',
    'siegen',
    'heben',
    'streifen',
    'drucken',
    'fordern',
    'greifen',
    'laufen',
    'zielen',
    'mühen',
    'drehen',
    'freuen',
    'scheiden',
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}

